# "did some update on reference page"
# Applies:
#  - Slide 11 (References): replace 4 raw URLs with descriptive titles
#  - Slide 6 (Complete history of background checks): crop + reposition picture
#  - Slide 7 (NICS Background Checks 2011-2017): reposition/format title,
#    crop + reposition picture
#  - Slide 9 (NIBRS Violent Crimes Data): reposition picture, nudge caption box

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 11 - "References:" - swap bare URLs for human-readable citation titles
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$refBody = $s11.Shapes.Item(2).TextFrame.TextRange

function Replace-Snippet($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -ge 0) {
        $sub = $textRange.Characters($idx + 1, $oldText.Length)
        $sub.Text = $newText
    }
}

Replace-Snippet $refBody "https://www.statista.com/statistics/249740/percentage-of-households-in-the-united-states-owning-a-firearm/" "Percentage of households in the United States owning one or more firearms from 1972 to 2018"
Replace-Snippet $refBody "https://www.statista.com/statistics/191962/participants-in-target-shooting-in-the-us-since-2006/" "Number of participants in target shooting in the United States from 2006 to 2017 (in millions)*"
Replace-Snippet $refBody "https://www.statista.com/statistics/195325/murder-victims-in-the-us-by-weapon-used/" "Number of murder victims in the United States in 2017, by weapon"
Replace-Snippet $refBody "https://injury.research.chop.edu/violence-prevention-initiative/types-violence-involving-youth/gun-violence/gun-violence-facts-and#.XTuao-hKhPY" "Gun Violence: Facts and Statistics"

# ---------------------------------------------------------------------------
# Slide 6 - "Complete history of background checks in the USA" - crop/move pic
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$pic6 = $s6.Shapes.Item(2)
$pic6.PictureFormat.CropBottom = 12.962
$pic6.Left = 92.27189976377953
$pic6.Top = 119.20007874015748
$pic6.Width = 752.8543707086615
$pic6.Height = 331.6737007874016

# ---------------------------------------------------------------------------
# Slide 7 - "NICS Background Checks 2011-2017" - resize/format title, crop/move pic
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

$title7 = $s7.Shapes.Item(1)
$title7.Left = 274.32648669291336
$title7.Top = 77.59223472440945
$title7.Width = 399.8446756692913
$title7.Height = 31.80586677165354
$title7.TextFrame.AutoSize = 2
$titleRange = $title7.TextFrame.TextRange
$titleRange.Font.Size = 18
$titleRange.Font.Bold = 1

$pic7 = $s7.Shapes.Item(2)
$pic7.PictureFormat.CropBottom = 13.7195
$pic7.Left = 88.9976577952756
$pic7.Top = 129.32033496062994
$pic7.Width = 756.1285826771654
$pic7.Height = 327.1455905511811

# ---------------------------------------------------------------------------
# Slide 9 - NIBRS Violent Crimes Data - widen picture, nudge callout box left
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$pic9 = $s9.Shapes.Item(1)
$pic9.Left = 130.01956755905513
$pic9.Width = 640.7712211023622

$callout9 = $s9.Shapes.Item(3)
$callout9.Left = 784.0723822047244
